$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before the old "Late" column (N), shifting
# Late / heading(Date) / Outstanding one column to the right (N->O, O->P, P->Q)
# to make room for the new "Variable Instalments" column.
$ws.Columns("N:N").Insert()

# The newly inserted column inherits the width of the column to its left (M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and select cell R6 on it
# (this also naturally clears the previous tab-selected state on "Transactions").
$ws.Activate()
$ws.Range("R6").Select()
